$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values scraped from coinranking.com at refresh time.
# Numeric-looking price cells get a leading apostrophe so Excel stores
# them as literal text (matching the sheet's existing inlineStr cells)
# instead of silently parsing them into floating-point numbers.

$ws.Range('D2').Value = '26.619.79'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.596.33'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range("D5").Value = "'211.67"
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range("D10").Value = "'19.51"
$ws.Range('E10').Value = '  -0.27%  '
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').Value = '1.819.85'
$ws.Range('E12').Value = '  +0.54%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'4.03"
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.556.41'
$ws.Range('E14').Value = '  -2.09%  '
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range("D16").Value = "'64.46"
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '26.606.48'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('E18').Value = '  +0.52%  '
$ws.Range("D19").Value = "'208.66"
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range("D21").Value = "'6.97"
$ws.Range('E21').Value = '  +3.80%  '
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('E23').Value = '  -1.84%  '
$ws.Range("D24").Value = "'8.89"
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range("D25").Value = "'145.23"
$ws.Range('E25').Value = '  -1.11%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range("D27").Value = "'7.14"
$ws.Range('E27').Value = '  -1.24%  '
$ws.Range('E28').Value = '  +0.62%  '
$ws.Range("D29").Value = "'15.27"
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  +0.76%  '
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range("D33").Value = "'0.656"
$ws.Range('E33').Value = '  -3.82%  '
$ws.Range('E34').Value = '  +1.10%  '
$ws.Range('D35').Value = '1.281.29'
$ws.Range('E35').Value = '  -1.91%  '
$ws.Range('E36').Value = '  +0.68%  '
$ws.Range('E37').Value = '  +1.05%  '
$ws.Range('E38').Value = '  -0.42%  '
$ws.Range('E39').Value = '  +1.84%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('E41').Value = '  +1.89%  '
$ws.Range("D42").Value = "'64.42"
$ws.Range('E42').Value = '  +2.83%  '
$ws.Range("D43").Value = "'0.786"
$ws.Range('E43').Value = '  -0.62%  '
$ws.Range('E44').Value = '  +1.29%  '
$ws.Range('D45').Value = '1.732.94'
$ws.Range('E45').Value = '  +0.55%  '
$ws.Range("D46").Value = "'0.910"
$ws.Range('E46').Value = '  +8.62%  '
$ws.Range("D47").Value = "'89.74"
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0105'
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = "'0.103"
$ws.Range('E50').Value = '  +5.15%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'0.0507"
$ws.Range('E51').Value = '  +0.54%  '
